$wb = $excel.ActiveWorkbook
$renters = $wb.Worksheets.Item("renters")
$owners = $wb.Worksheets.Item("owners")

# ---------------------------------------------------------------------
# owners sheet: a couple of owner savings amounts grew, two damage
# states got reset to "None", and M8 (home value) became a formula
# (K8*279) instead of a hard-coded number, which ripples into F8.
# ---------------------------------------------------------------------
$owners.Range("C2").Value2 = 100000
$owners.Range("N3").Value2 = "None"
$owners.Range("C4").Value2 = 50000
$owners.Range("N5").Value2 = "None"
$owners.Range("M8").Formula = "=K8*279"
$owners.Range("N13").Select()

# ---------------------------------------------------------------------
# renters sheet: two newly-found rental homes get appended for tenants
# who were previously in find_home (now rent_home / entities.HomeOwner()).
# Cell-set order matters because it controls shared-string append order.
# ---------------------------------------------------------------------

# Row 10 - For Alfred
$renters.Range("B10").Value2 = "341 Where St"
$renters.Range("P10").Value2 = "Gal"
$renters.Range("A10").Value2 = "For Alfred"
$renters.Range("C10").Value2 = "Mobile Home"
$renters.Range("D10").Value2 = "Rental"
$renters.Range("E10").Value2 = 60000
$renters.Range("F10").Value2 = 1000
$renters.Range("G10").Value2 = 700
$renters.Range("H10").Value2 = 0
$renters.Range("I10").Value2 = 1
$renters.Range("I10").NumberFormat = "#,##0_);[Red](#,##0)"
$renters.Range("J10").Value2 = 1
$renters.Range("K10").Value2 = 1
$renters.Range("L10").Value2 = 1100
$renters.Range("M10").Value2 = 1920
$renters.Range("N10").Formula = "=L10*279"
$renters.Range("O10").Value2 = "None"
$renters.Range("Q10").Value2 = 30000
$renters.Range("R10").Value2 = 0
$renters.Range("S10").Value2 = $true
$renters.Range("T10").Value2 = -90.294238000000007
$renters.Range("U10").Value2 = 43.224015000000001
$renters.Range("V10").Value2 = 700

# Row 11 - For Selena
$renters.Range("A11").Value2 = "For Selena"
$renters.Range("B11").Value2 = "9900 Nowhere St"
$renters.Range("B11").Font.Color = 0
$renters.Range("C11").Value2 = "Single Family Dwelling"
$renters.Range("D11").Value2 = "Rental"
$renters.Range("D11").Font.Color = 0
$renters.Range("E11").Value2 = 60000
$renters.Range("E11").Font.Color = 0
$renters.Range("F11").Value2 = 1000
$renters.Range("F11").Font.Color = 0
$renters.Range("G11").Value2 = 700
$renters.Range("G11").Font.Color = 0
$renters.Range("H11").Value2 = 0
$renters.Range("H11").Font.Color = 0
$renters.Range("I11").Value2 = 1
$renters.Range("I11").Font.Color = 0
$renters.Range("I11").NumberFormat = "#,##0_);[Red](#,##0)"
$renters.Range("J11").Value2 = 2
$renters.Range("K11").Value2 = 1
$renters.Range("L11").Value2 = 750
$renters.Range("M11").Value2 = 1960
$renters.Range("N11").Formula = "=L11*279"
$renters.Range("O11").Value2 = "None"
$renters.Range("P11").Value2 = "Blake"
$renters.Range("Q11").Value2 = 30000
$renters.Range("Q11").Font.Color = 0
$renters.Range("R11").Value2 = 0
$renters.Range("R11").Font.Color = 0
$renters.Range("S11").Value2 = $true
$renters.Range("S11").Font.Color = 0
$renters.Range("T11").Value2 = -90.294238000000007
$renters.Range("T11").Font.Color = 0
$renters.Range("U11").Value2 = 43.224015000000001
$renters.Range("U11").Font.Color = 0
$renters.Range("V11").Value2 = 700
$renters.Range("V11").Font.Color = 0

$renters.Range("H22").Select()
